$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "add"
$ws.Range("A2").Value = "minus"
$ws.Range("A3").Value = "plus"

$ws.Range("A3").Select()
